# Refresh cryptos list (coinranking.com scrape) with the latest prices,
# 1h volume deltas, and a handful of re-ranked coins.
#
# Column D ("Price") is stored as literal display text (e.g. "3.592.86" is
# a euro-style thousands separator, not a number). Excel's COM layer auto-
# coerces a plain numeric-looking string assigned via .Value into a real
# Number (dropping the trailing zeros we need, e.g. "78.00" -> 78). Prefix
# those assignments with a literal leading apostrophe (') -- exactly what
# typing into the Excel UI does -- to force a Text cell that keeps the
# string byte-for-byte. Non-ambiguous strings (extra "." thousands seps,
# URLs, names, the "  +/-x.xx%  " volume column) don't need it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.286.73'
$ws.Range("E2").Value = '  -2.07%  '

$ws.Range("D3").Value = '3.603.58'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''625.84'
$ws.Range("E5").Value = '  -7.04%  '

$ws.Range("D6").Value = '''156.78'
$ws.Range("E6").Value = '  -2.66%  '

$ws.Range("D7").Value = '3.601.92'
$ws.Range("E7").Value = '  -2.34%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = '''0.489'
$ws.Range("E9").Value = '  -2.09%  '

$ws.Range("E10").Value = '  -2.75%  '

$ws.Range("D11").Value = '''7.02'
$ws.Range("E11").Value = '  -1.10%  '

$ws.Range("D12").Value = '''0.435'
$ws.Range("E12").Value = '  -1.90%  '

$ws.Range("D13").Value = '''0.0000226'
$ws.Range("E13").Value = '  -3.46%  '

$ws.Range("D14").Value = '4.213.40'
$ws.Range("E14").Value = '  -2.43%  '

$ws.Range("D15").Value = '''32.06'
$ws.Range("E15").Value = '  -3.47%  '

$ws.Range("D16").Value = '3.591.52'
$ws.Range("E16").Value = '  -2.99%  '

$ws.Range("D17").Value = '68.236.88'
$ws.Range("E17").Value = '  -2.09%  '

$ws.Range("D18").Value = '''0.118'
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").Value = '''6.43'
$ws.Range("E19").Value = '  -1.25%  '

$ws.Range("D20").Value = '''15.65'
$ws.Range("E20").Value = '  -3.24%  '

$ws.Range("D21").Value = '''9.92'
$ws.Range("E21").Value = '  +1.36%  '

$ws.Range("D22").Value = '''458.29'
$ws.Range("E22").Value = '  -2.81%  '

$ws.Range("D23").Value = '''0.643'
$ws.Range("E23").Value = '  -0.93%  '

$ws.Range("D24").Value = '''78.00'
$ws.Range("E24").Value = '  -2.55%  '

$ws.Range("D25").Value = '3.741.72'
$ws.Range("E25").Value = '  -2.57%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").Value = '''10.77'
$ws.Range("E27").Value = '  -1.86%  '

$ws.Range("E28").Value = '  -8.58%  '

$ws.Range("D29").Value = '''8.45'
$ws.Range("E29").Value = '  -7.32%  '

$ws.Range("D30").Value = '''2.60'
$ws.Range("E30").Value = '  -3.76%  '

$ws.Range("D31").Value = '''1.64'
$ws.Range("E31").Value = '  -4.73%  '

$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''1.93'
$ws.Range("E33").Value = '  -4.96%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''26.13'
$ws.Range("E34").Value = '  -2.74%  '

$ws.Range("D35").Value = '''0.160'
$ws.Range("E35").Value = '  -5.00%  '

$ws.Range("D36").Value = '3.600.56'
$ws.Range("E36").Value = '  -2.40%  '

$ws.Range("D37").Value = '''6.24'
$ws.Range("E37").Value = '  -4.33%  '

$ws.Range("D38").Value = '''8.20'
$ws.Range("E38").Value = '  -3.54%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''177.42'
$ws.Range("E40").Value = '  +0.81%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.998'
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").Value = '''5.66'
$ws.Range("E42").Value = '  -7.93%  '

$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").Value = '''0.0887'
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '''2.16'
$ws.Range("E44").Value = '  -5.34%  '

$ws.Range("D45").Value = '''0.908'
$ws.Range("E45").Value = '  -2.92%  '

$ws.Range("D46").Value = '''29.08'
$ws.Range("E46").Value = '  +4.17%  '

$ws.Range("D47").Value = '''46.04'
$ws.Range("E47").Value = '  -2.01%  '

$ws.Range("D48").Value = '''2.61'
$ws.Range("E48").Value = '  -5.50%  '

$ws.Range("D49").Value = '''7.74'
$ws.Range("E49").Value = '  -2.15%  '

$ws.Range("E50").Value = '  -5.43%  '

$ws.Range("D51").Value = '''1.20'
$ws.Range("E51").Value = '  -7.36%  '
